$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.434.02'
$ws.Cells.Item(2, 5).Value = '  +4.16%  '

$ws.Cells.Item(3, 4).Value = '1.596.89'
$ws.Cells.Item(3, 5).Value = '  +1.94%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '214.82'
$ws.Cells.Item(5, 5).Value = '  +2.10%  '

$ws.Cells.Item(6, 5).Value = '  +1.76%  '

$ws.Cells.Item(7, 5).Value = '  -0.07%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '24.04'
$ws.Cells.Item(8, 5).Value = '  +8.89%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.252'
$ws.Cells.Item(9, 5).Value = '  +1.29%  '

$ws.Cells.Item(10, 5).Value = '  +0.84%  '

$ws.Cells.Item(11, 5).Value = '  +2.33%  '

$ws.Cells.Item(12, 4).Value = '1.821.73'
$ws.Cells.Item(12, 5).Value = '  +1.79%  '

$ws.Cells.Item(13, 4).Value = '1.600.84'
$ws.Cells.Item(13, 5).Value = '  +2.28%  '

$ws.Cells.Item(14, 5).Value = '  +0.84%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.534'
$ws.Cells.Item(15, 5).Value = '  +3.13%  '

$ws.Cells.Item(16, 4).Value = '28.466.31'
$ws.Cells.Item(16, 5).Value = '  +4.40%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '63.33'
$ws.Cells.Item(17, 5).Value = '  +2.34%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '233.02'
$ws.Cells.Item(18, 5).Value = '  +7.32%  '

$ws.Cells.Item(20, 5).Value = '  +0.87%  '

$ws.Cells.Item(21, 5).Value = '  -0.16%  '

$ws.Cells.Item(22, 5).Value = '  -0.13%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '9.43'
$ws.Cells.Item(23, 5).Value = '  +2.43%  '

$ws.Cells.Item(24, 5).Value = '  +1.09%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '152.70'
$ws.Cells.Item(25, 5).Value = '  -0.09%  '

$ws.Cells.Item(26, 5).Value = '  +2.03%  '

$ws.Cells.Item(27, 5).Value = '  -0.05%  '

$ws.Cells.Item(28, 5).Value = '  +1.18%  '

$ws.Cells.Item(29, 5).Value = '  -0.07%  '

$ws.Cells.Item(30, 5).Value = '  +0.93%  '

$ws.Cells.Item(31, 5).Value = '  +1.29%  '

$ws.Cells.Item(32, 5).Value = '  +0.41%  '

$ws.Cells.Item(33, 5).Value = '  +0.77%  '

$ws.Cells.Item(34, 4).Value = '1.420.80'
$ws.Cells.Item(34, 5).Value = '  -0.91%  '

$ws.Cells.Item(35, 5).Value = '  -0.68%  '

$ws.Cells.Item(36, 5).Value = '  -3.63%  '

$ws.Cells.Item(37, 5).Value = '  -0.01%  '

$ws.Cells.Item(38, 5).Value = '  +0.64%  '

$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.547'
$ws.Cells.Item(39, 5).Value = '  +2.72%  '

$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.53'
$ws.Cells.Item(40, 5).Value = '  +8.26%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.821'
$ws.Cells.Item(41, 5).Value = '  +1.71%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '5.76'
$ws.Cells.Item(42, 5).Value = '  -3.03%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.981'
$ws.Cells.Item(44, 5).Value = '  -1.57%  '

$ws.Cells.Item(45, 5).Value = '  +6.27%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '64.88'
$ws.Cells.Item(46, 5).Value = '  +0.62%  '

$ws.Cells.Item(47, 4).Value = '1.734.40'
$ws.Cells.Item(47, 5).Value = '  +1.84%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '87.63'
$ws.Cells.Item(48, 5).Value = '  +1.80%  '

$ws.Cells.Item(49, 5).Value = '  +0.29%  '

$ws.Cells.Item(50, 5).Value = '  +5.19%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.0524'
$ws.Cells.Item(51, 5).Value = '  -0.19%  '
